$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = "Jane"
$ws.Range("P2").Select()
